$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 94 (existing rows 94:106 shift down to 96:108)
$ws.Range("A94:A95").EntireRow.Insert()

# Common (constant across rows) columns for this data block
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria = "Ciruela"
$unidad = "$/bandeja 18 kilos granel"
$origen = "Región de O'Higgins"

# New row 94: Angeleno, Especial
$r = 94
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45013
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Angeleno"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 70
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = 18

# New row 95: Angeleno, Primera
$r = 95
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45013
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Angeleno"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 60
$ws.Cells.Item($r, 14).Value = 11000
$ws.Cells.Item($r, 15).Value = 11000
$ws.Cells.Item($r, 16).Value = 11000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 611
$ws.Cells.Item($r, 20).Value = 18
